$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number (e.g. "141.97").
# These must be forced to Text format first, otherwise Excel auto-converts the
# assigned string into a numeric value, which would change the cell type from
# the original inline string ("t=inlineStr" / shared-string text) to a number.
$textForcedCells = @(
    "D5"
    "D6"
    "D12"
    "D16"
    "D22"
    "D25"
    "D28"
    "D30"
    "D33"
    "D34"
    "D35"
    "D36"
    "D37"
    "D41"
    "D42"
    "D43"
    "D46"
    "D47"
    "D48"
    "D50"
    "D51"
)
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated Price (column D) values that required the text coercion above.
$ws.Range("D5").Value = "568.58"
$ws.Range("D6").Value = "141.97"
$ws.Range("D12").Value = "0.368"
$ws.Range("D16").Value = "23.37"
$ws.Range("D22").Value = "6.99"
$ws.Range("D25").Value = "63.14"
$ws.Range("D28").Value = "7.68"
$ws.Range("D30").Value = "1.83"
$ws.Range("D33").Value = "160.96"
$ws.Range("D34").Value = "19.47"
$ws.Range("D35").Value = "4.23"
$ws.Range("D36").Value = "0.968"
$ws.Range("D37").Value = "1.21"
$ws.Range("D41").Value = "0.849"
$ws.Range("D42").Value = "294.23"
$ws.Range("D43").Value = "137.45"
$ws.Range("D46").Value = "0.0982"
$ws.Range("D47").Value = "19.69"
$ws.Range("D48").Value = "0.0545"
$ws.Range("D50").Value = "19.88"
$ws.Range("D51").Value = "10.71"

# Restore the default "Normal" style on those cells so no stray number format
# is left behind on the cell itself (the underlying value remains text).
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining Price/Volume(1h) text updates that are already unambiguous text
# (contain extra separators, percent signs, or padding spaces) and therefore
# do not get reinterpreted as numbers by Excel.
$ws.Range("D2").Value = "60.474.57"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "2.597.97"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "2.618.94"
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("E12").Value = "  +3.99%  "
$ws.Range("E13").Value = "  -6.90%  "
$ws.Range("D14").Value = "3.061.52"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "60.433.61"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").Value = "2.609.38"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("E19").Value = "  +9.04%  "
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("E22").Value = "  +9.53%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("E24").Value = "  +14.02%  "
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("E28").Value = "  +4.93%  "
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("E30").Value = "  +9.24%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("E35").Value = "  +4.77%  "
$ws.Range("E36").Value = "  +10.25%  "
$ws.Range("E37").Value = "  +4.42%  "
$ws.Range("E38").Value = "  +8.38%  "
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("E40").Value = "  +3.89%  "
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  +4.01%  "
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("E47").Value = "  +3.65%  "
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("E50").Value = "  +6.53%  "
$ws.Range("E51").Value = "  +0.54%  "
